$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.421.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E7").Value = "  -6.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3835"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.14"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -10.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07819"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.013"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.61%  "
$ws.Range("E12").Value = "  -3.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.856.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.837"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.092"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.17%  "
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "85.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.33%  "
$ws.Range("E18").Value = "  -4.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06496"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.28%  "
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.467"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.406.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.95%  "
$ws.Range("E24").Value = "  -7.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.284"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.059.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("E28").Value = "  -4.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.466"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.026"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.492"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09333"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9257"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.618"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.204"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.22%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02216"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.13%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.214"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05942"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.273"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.21%  "
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5887"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1846"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.261"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5619"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.04%  "
$ws.Range("E47").Value = "  -6.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.354"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.904"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06839"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("E51").Value = "  -0.83%  "
